$wb = $excel.ActiveWorkbook

# --- FilterRecord sheet: add two new filter record values ---
$wsFilterRecord = $wb.Worksheets.Item("FilterRecord")
$wsFilterRecord.Range("A4").Value = "CVAS - IP Valuation"
$wsFilterRecord.Range("A5").Value = "Lender Education"
$wsFilterRecord.Columns.Item(1).ColumnWidth = 37.2

# --- ReportOption sheet: add two more "Opportunities" rows ---
$wsReportOption = $wb.Worksheets.Item("ReportOption")
$wsReportOption.Range("A4").Value = "Opportunities"
$wsReportOption.Range("A5").Value = "Opportunities"

# --- Filter sheet: add two more "Job Type" rows ---
$wsFilter = $wb.Worksheets.Item("Filter")
$wsFilter.Range("A4").Value = "Job Type"
$wsFilter.Range("A5").Value = "Job Type"

# --- Update sheet selections / active sheet to match the saved view state ---
$wsReportOption.Activate()
$wsReportOption.Range("A5").Select()

$wsFilter.Activate()
$wsFilter.Range("A5").Select()

$wsFilterRecord.Activate()
$wsFilterRecord.Range("C14").Select()
